$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the six oldest year rows (2000, 2005, 2006, 2007, 2008, 2009).
# This shifts the remaining rows (2010-2013 data) up so they become rows 2-5.
$ws.Range("A2:E7").EntireRow.Delete()
